$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "(1,2, p = 11, c = 400)"
$ws.Range("C2").Value = 11

$ws.Range("B3").Value = "(1,4, p = 7, c = 133)"
$ws.Range("C3").Value = 7

$ws.Range("B4").Value = "(2,3, p = 10, c = 700)"
$ws.Range("C4").Value = 10

$ws.Range("B5").Value = "(2,5, p = 16, c = 133)"
$ws.Range("C5").Value = 16

$ws.Range("B6").Value = "(3,6, p = 6, c = 400)"
$ws.Range("C6").Value = 6

$ws.Range("B7").Value = "(4,5, p = 3, c = 700)"
$ws.Range("C7").Value = 3

$ws.Range("B8").Value = "(4,7, p = 10, c = 100)"
$ws.Range("C8").Value = 10

$ws.Range("B9").Value = "(5,6, p = 10, c = 600)"
$ws.Range("C9").Value = 10

$ws.Range("B10").Value = "(5,8, p = 15, c = 700)"
$ws.Range("C10").Value = 15
